# Auto-generated: update Leve profit/price cells per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 22066.6
$ws.Range("I21").Value = 20833.25
$ws.Range("J21").Value = 27000
$ws.Range("K21").Value = 20833.25
$ws.Range("L21").Value = 27000
$ws.Range("M21").Value = -20365.25
$ws.Range("N21").Value = -27936
$ws.Range("H23").Value = 22066.6
$ws.Range("I23").Value = 20833.25
$ws.Range("J23").Value = 27000
$ws.Range("K23").Value = 20833.25
$ws.Range("L23").Value = 27000
$ws.Range("M23").Value = -20599.25
$ws.Range("N23").Value = -27468
$ws.Range("H33").Value = 237.46666
$ws.Range("I33").Value = 108.888885
$ws.Range("J33").Value = 430.33334
$ws.Range("K33").Value = 108.888885
$ws.Range("L33").Value = 430.33334
$ws.Range("M33").Value = 120.111115
$ws.Range("N33").Value = -888.33334
$ws.Range("H129").Value = 915.4423
$ws.Range("J129").Value = 1037.6904
$ws.Range("L129").Value = 3113.0712
$ws.Range("N129").Value = -13113.0712
$ws.Range("I137").Value = 90911096
$ws.Range("J137").Value = 1810.3
$ws.Range("K137").Value = 272733288
$ws.Range("L137").Value = 5430.9
$ws.Range("M137").Value = -272730738
$ws.Range("N137").Value = -10530.9
$ws.Range("H141").Value = 1967.45
$ws.Range("J141").Value = 4401.6665
$ws.Range("L141").Value = 13204.9995
$ws.Range("N141").Value = -23564.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3851.081
$ws.Range("I32").Value = 2703.3103
$ws.Range("K32").Value = 2703.3103
$ws.Range("M32").Value = -2416.3103
$ws.Range("H45").Value = 1277.7778
$ws.Range("I45").Value = 1266.6666
$ws.Range("K45").Value = 1266.6666
$ws.Range("M45").Value = -889.6666
$ws.Range("H61").Value = 2533.9697
$ws.Range("I61").Value = 1803.4783
$ws.Range("J61").Value = 4214.1
$ws.Range("K61").Value = 1803.4783
$ws.Range("L61").Value = 4214.1
$ws.Range("M61").Value = -1591.4783
$ws.Range("N61").Value = -4638.1
$ws.Range("H122").Value = 1676
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150
$ws.Range("H123").Value = 32952
$ws.Range("J123").Value = 32952
$ws.Range("L123").Value = 32952
$ws.Range("N123").Value = -42752
$ws.Range("H136").Value = 2533.9697
$ws.Range("I136").Value = 1803.4783
$ws.Range("J136").Value = 4214.1
$ws.Range("K136").Value = 5410.4349
$ws.Range("L136").Value = 12642.3
$ws.Range("M136").Value = -2860.4349
$ws.Range("N136").Value = -17742.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 908.8570999999999
$ws.Range("I5").Value = 531.4
$ws.Range("K5").Value = 1594.2
$ws.Range("M5").Value = -1482.2
$ws.Range("H23").Value = 682.5
$ws.Range("I23").Value = 1158.7778
$ws.Range("J23").Value = 206.22223
$ws.Range("K23").Value = 3476.3334
$ws.Range("L23").Value = 618.66669
$ws.Range("M23").Value = -3241.3334
$ws.Range("N23").Value = -1088.66669
$ws.Range("H113").Value = 12195892
$ws.Range("J113").Value = 13514304
$ws.Range("L113").Value = 40542912
$ws.Range("N113").Value = -40547252
$ws.Range("H122").Value = 1171.2
$ws.Range("I122").Value = 804
$ws.Range("J122").Value = 1212
$ws.Range("K122").Value = 7236
$ws.Range("L122").Value = 10908
$ws.Range("M122").Value = -4786
$ws.Range("N122").Value = -15808
$ws.Range("H131").Value = 2734.6448
$ws.Range("I131").Value = 365
$ws.Range("J131").Value = 2798.6892
$ws.Range("K131").Value = 1095
$ws.Range("L131").Value = 8396.067599999998
$ws.Range("M131").Value = 3945
$ws.Range("N131").Value = -18476.0676
$ws.Range("H135").Value = 908.8570999999999
$ws.Range("I135").Value = 531.4
$ws.Range("K135").Value = 4782.599999999999
$ws.Range("M135").Value = -2247.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25000
$ws.Range("I57").Value = 25000
$ws.Range("K57").Value = 25000
$ws.Range("M57").Value = -24180
$ws.Range("H70").Value = 6762.5264
$ws.Range("I70").Value = 7422.154
$ws.Range("J70").Value = 5333.3335
$ws.Range("K70").Value = 7422.154
$ws.Range("L70").Value = 5333.3335
$ws.Range("M70").Value = -7152.154
$ws.Range("N70").Value = -5873.3335
$ws.Range("H73").Value = 6762.5264
$ws.Range("I73").Value = 7422.154
$ws.Range("J73").Value = 5333.3335
$ws.Range("K73").Value = 7422.154
$ws.Range("L73").Value = 5333.3335
$ws.Range("M73").Value = -6486.154
$ws.Range("N73").Value = -7205.3335
$ws.Range("H122").Value = 2779262.8
$ws.Range("I122").Value = 11111111
$ws.Range("J122").Value = 1980
$ws.Range("K122").Value = 33333333
$ws.Range("L122").Value = 5940
$ws.Range("M122").Value = -33330883
$ws.Range("N122").Value = -10840
$ws.Range("H132").Value = 2415.2144
$ws.Range("I132").Value = 2039.4762
$ws.Range("J132").Value = 3542.4285
$ws.Range("K132").Value = 6118.4286
$ws.Range("L132").Value = 10627.2855
$ws.Range("M132").Value = -3588.4286
$ws.Range("N132").Value = -15687.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3067
$ws.Range("I7").Value = 2067.7144
$ws.Range("J7").Value = 3533.3333
$ws.Range("K7").Value = 2067.7144
$ws.Range("L7").Value = 3533.3333
$ws.Range("M7").Value = -1955.7144
$ws.Range("N7").Value = -3757.3333
$ws.Range("H40").Value = 3036.842
$ws.Range("I40").Value = 1450
$ws.Range("J40").Value = 3223.5293
$ws.Range("K40").Value = 1450
$ws.Range("L40").Value = 3223.5293
$ws.Range("M40").Value = -1314
$ws.Range("N40").Value = -3495.5293
$ws.Range("H61").Value = 4128.8887
$ws.Range("I61").Value = 5812.222
$ws.Range("J61").Value = 2445.5557
$ws.Range("K61").Value = 5812.222
$ws.Range("L61").Value = 2445.5557
$ws.Range("M61").Value = -5610.222
$ws.Range("N61").Value = -2849.5557
$ws.Range("H113").Value = 4128.8887
$ws.Range("I113").Value = 5812.222
$ws.Range("J113").Value = 2445.5557
$ws.Range("K113").Value = 5812.222
$ws.Range("L113").Value = 2445.5557
$ws.Range("M113").Value = -3642.222
$ws.Range("N113").Value = -6785.5557
$ws.Range("H122").Value = 3731.7896
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3731.7896
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11195.3688
$ws.Range("N122").Value = -16095.3688
$ws.Range("H126").Value = 3067
$ws.Range("I126").Value = 2067.7144
$ws.Range("J126").Value = 3533.3333
$ws.Range("K126").Value = 6203.1432
$ws.Range("L126").Value = 10599.9999
$ws.Range("M126").Value = -3733.1432
$ws.Range("N126").Value = -15539.9999
$ws.Range("H135").Value = 34800
$ws.Range("J135").Value = 34800
$ws.Range("L135").Value = 34800
$ws.Range("N135").Value = -44940
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 168317.33
$ws.Range("I122").Value = 502252
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 1506756
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -1504306
$ws.Range("N122").Value = -8950
$ws.Range("H123").Value = 47255
$ws.Range("J123").Value = 47255
$ws.Range("L123").Value = 47255
$ws.Range("N123").Value = -57055
